# Advance Feature file task updates
#
# - On the "profile" sheet: the "AvailabilityHours" answer in row 3 moves
#   from "More than 30hours a week" to "Less than 30hours a week", and a
#   new row 4 is added holding a plain numeric value (2) in column B.
# - The active cell on the "profile" sheet becomes B4.
# - "Search-Skills" stays the active/selected worksheet in the workbook,
#   exactly as it was before the edit.

$wb = $excel.ActiveWorkbook

$profile = $wb.Worksheets.Item("profile")
$searchSkills = $wb.Worksheets.Item("Search-Skills")

# Update the availability-hours text for the "Part Time" row.
$profile.Range("B3").Value = "Less than 30hours a week"

# New row with a numeric value.
$profile.Range("B4").Value = 2

# Move the selection on the profile sheet to the newly added cell. This
# requires briefly activating the sheet so the selection can be set.
$profile.Activate()
$profile.Range("B4").Select()

# Restore "Search-Skills" as the active sheet/tab, matching the original
# workbook state (only the profile sheet's own selection should change).
$searchSkills.Activate()
